$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 88
$wsExpo.Range("F3").Value = 357
$wsExpo.Range("F4").Value = 4738
$wsExpo.Range("F6").Value = 478

# Sheet "全部类型" (all types) - same events repeated, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 88
$wsAll.Range("F3").Value = 357
$wsAll.Range("F4").Value = 4738
$wsAll.Range("F8").Value = 478
